$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6667055.5
$ws.Range("J33").Value = 1179.6666
$ws.Range("L33").Value = 1179.6666
$ws.Range("N33").Value = -1637.6666
$ws.Range("H45").Value = 15019
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 15019
$ws.Range("K45").Value = 0
$ws.Range("L45").ClearContents()
$ws.Range("M45").Value = 45057
$ws.Range("N45").Value = -45441
$ws.Range("H46").Value = 3000
$ws.Range("J46").Value = 1000
$ws.Range("L46").Value = 3000
$ws.Range("N46").Value = -3238
$ws.Range("H60").Value = 3000
$ws.Range("J60").Value = 1000
$ws.Range("L60").Value = 3000
$ws.Range("N60").Value = -3968
$ws.Range("H109").Value = 47341
$ws.Range("J109").Value = 47341
$ws.Range("L109").Value = 47341
$ws.Range("N109").Value = -50115
$ws.Range("H112").Value = 78586.38
$ws.Range("J112").Value = 102007.9
$ws.Range("L112").Value = 306023.7
$ws.Range("N112").Value = -308239.7
$ws.Range("H137").Value = 2098.25
$ws.Range("I137").Value = 1975.3529
$ws.Range("K137").Value = 5926.0587
$ws.Range("M137").Value = -3376.0587

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 35000
$ws.Range("J24").Value = 35000
$ws.Range("L24").Value = 35000
$ws.Range("N24").Value = -35748
$ws.Range("H76").Value = 25000
$ws.Range("J76").Value = 25000
$ws.Range("L76").Value = 25000
$ws.Range("N76").Value = -25676
$ws.Range("H79").Value = 25000
$ws.Range("J79").Value = 25000
$ws.Range("L79").Value = 25000
$ws.Range("N79").Value = -27340
$ws.Range("H92").Value = 29000
$ws.Range("J92").Value = 29000
$ws.Range("L92").Value = 29000
$ws.Range("N92").Value = -33992
$ws.Range("H97").Value = 1239.3125
$ws.Range("I97").Value = 1239.3125
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1239.3125
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -743.3125
$ws.Range("H100").Value = 35000
$ws.Range("J100").Value = 35000
$ws.Range("L100").Value = 35000
$ws.Range("N100").Value = -37164
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").ClearContents()
$ws.Range("N107").Value = 0
$ws.Range("H122").Value = 1729.75
$ws.Range("I122").Value = 1729.75
$ws.Range("K122").Value = 5189.25
$ws.Range("M122").Value = -2739.25
$ws.Range("H132").Value = 3182.1785
$ws.Range("I132").Value = 2008.3684
$ws.Range("J132").Value = 5660.222
$ws.Range("K132").Value = 6025.1052
$ws.Range("L132").Value = 16980.666
$ws.Range("M132").Value = -3495.1052
$ws.Range("N132").Value = -22040.666

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1703.0698
$ws.Range("I94").Value = 1509.1316
$ws.Range("K94").Value = 1509.1316
$ws.Range("M94").Value = -1058.1316
$ws.Range("H107").Value = 1796.5883
$ws.Range("I107").Value = 1735.75
$ws.Range("J107").Value = 1942.6
$ws.Range("K107").Value = 1735.75
$ws.Range("L107").Value = 1942.6
$ws.Range("M107").Value = 184.25
$ws.Range("N107").Value = -5782.6

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2904.818
$ws.Range("J15").Value = 4429.7144
$ws.Range("L15").Value = 4429.7144
$ws.Range("N15").Value = -4769.7144
$ws.Range("H31").Value = 27438.676
$ws.Range("J31").Value = 10761.333
$ws.Range("L31").Value = 10761.333
$ws.Range("N31").Value = -11351.333
$ws.Range("H34").Value = 27438.676
$ws.Range("J34").Value = 10761.333
$ws.Range("L34").Value = 10761.333
$ws.Range("N34").Value = -11165.333
$ws.Range("H58").Value = 2136.95
$ws.Range("I58").Value = 2239.9714
$ws.Range("K58").Value = 2239.9714
$ws.Range("M58").Value = -2036.9714
$ws.Range("H134").Value = 8709.974
$ws.Range("I134").Value = 6084.4443
$ws.Range("J134").Value = 15154.454
$ws.Range("K134").Value = 18253.3329
$ws.Range("L134").Value = 45463.362
$ws.Range("M134").Value = -15718.3329
$ws.Range("N134").Value = -50533.362
$ws.Range("H136").Value = 2136.95
$ws.Range("I136").Value = 2239.9714
$ws.Range("K136").Value = 6719.914199999999
$ws.Range("M136").Value = -4169.914199999999

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 12408.667
$ws.Range("I115").Value = 7983.6
$ws.Range("J115").Value = 17940
$ws.Range("K115").Value = 23950.8
$ws.Range("L115").Value = 53820
$ws.Range("M115").Value = -22775.8
$ws.Range("N115").Value = -56170
$ws.Range("H116").Value = 5562.278
$ws.Range("J116").Value = 5574.8667
$ws.Range("L116").Value = 16724.6001
$ws.Range("N116").Value = -23608.6001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7501250.5
$ws.Range("I11").Value = 12000600
$ws.Range("J11").Value = 2334
$ws.Range("K11").Value = 12000600
$ws.Range("L11").Value = 2334
$ws.Range("M11").Value = -12000461
$ws.Range("N11").Value = -2612
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H14").Value = 2601081.8
$ws.Range("I14").Value = 3714916.2
$ws.Range("J14").Value = 2134.6667
$ws.Range("K14").Value = 3714916.2
$ws.Range("L14").Value = 2134.6667
$ws.Range("M14").Value = -3714748.2
$ws.Range("N14").Value = -2470.6667

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = 0
$ws.Range("H17").Value = 25250.5
$ws.Range("I17").Value = 25250.5
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 25250.5
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -25080.5
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").ClearContents()
$ws.Range("N38").Value = 0
$ws.Range("H40").Value = 4702.086
$ws.Range("I40").Value = 4186.8
$ws.Range("J40").Value = 5990.3
$ws.Range("K40").Value = 4186.8
$ws.Range("L40").Value = 5990.3
$ws.Range("M40").Value = -4050.8
$ws.Range("N40").Value = -6262.3
$ws.Range("H110").Value = 70000
$ws.Range("J110").Value = 70000
$ws.Range("L110").Value = 70000
$ws.Range("N110").Value = -78180

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 682.5263
$ws.Range("I113").Value = 618.8
$ws.Range("K113").Value = 1856.4
$ws.Range("M113").Value = 313.6000000000001
$ws.Range("H132").Value = 3185.0435
$ws.Range("I132").Value = 3052.4167
$ws.Range("J132").Value = 3662.5
$ws.Range("K132").Value = 9157.250100000001
$ws.Range("L132").Value = 10987.5
$ws.Range("M132").Value = -6627.250100000001
$ws.Range("N132").Value = -16047.5
$ws.Range("H136").Value = 1751
$ws.Range("I136").Value = 1692.25
$ws.Range("J136").Value = 2080
$ws.Range("K136").Value = 5076.75
$ws.Range("L136").Value = 6240
$ws.Range("M136").Value = -2526.75
$ws.Range("N136").Value = -11340
$ws.Range("H137").Value = 96047.336
$ws.Range("J137").Value = 96047.336
$ws.Range("L137").Value = 96047.336
$ws.Range("N137").Value = -106247.336
